# Update gh-pages to output generated at 456a3b4
# Refresh "想去人数" (wanted-to-go count) figures and flip one sold-out
# show's "最低票价" (min ticket price) to the literal "已售罄" label,
# across the 展览 / 演出 / 全部类型 sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (Exhibitions) ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 164
$ws.Range("F6").Value  = 1374
$ws.Range("F10").Value = 486
$ws.Range("F11").Value = 843
$ws.Range("F12").Value = 546
$ws.Range("F13").Value = 761
$ws.Range("F14").Value = 332
$ws.Range("F15").Value = 506
$ws.Range("F17").Value = 1074
$ws.Range("F19").Value = 307
$ws.Range("F22").Value = 255
$ws.Range("F23").Value = 35
$ws.Range("F25").Value = 504
$ws.Range("F26").Value = 476
$ws.Range("F28").Value = 371

# ---- Sheet "演出" (Performances) ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("G2").Value  = "已售罄"
$ws.Range("F4").Value  = 58
$ws.Range("F5").Value  = 51
$ws.Range("F6").Value  = 299
$ws.Range("F10").Value = 165
$ws.Range("F11").Value = 167

# ---- Sheet "全部类型" (All types) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 164
$ws.Range("F7").Value  = 1374
$ws.Range("G8").Value  = "已售罄"
$ws.Range("F12").Value = 58
$ws.Range("F14").Value = 51
$ws.Range("F15").Value = 299
$ws.Range("F16").Value = 486
$ws.Range("F17").Value = 843
$ws.Range("F18").Value = 546
$ws.Range("F19").Value = 761
$ws.Range("F20").Value = 332
$ws.Range("F21").Value = 506
$ws.Range("F23").Value = 1074
$ws.Range("F27").Value = 307
$ws.Range("F31").Value = 165
$ws.Range("F32").Value = 255
$ws.Range("F33").Value = 35
$ws.Range("F35").Value = 167
$ws.Range("F37").Value = 504
$ws.Range("F40").Value = 476
$ws.Range("F42").Value = 371
